$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, pushing existing rows 59..141 down to 60..142
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly record
$ws.Cells.Item(59, 1).Value = 8
$ws.Cells.Item(59, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44679
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = 100112044
$ws.Cells.Item(59, 7).Value = "Perejil"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 3000
$ws.Cells.Item(59, 11).Value = 2500
$ws.Cells.Item(59, 12).Value = 3000
$ws.Cells.Item(59, 13).Value = 2750
$ws.Cells.Item(59, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(59, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(59, 16).Value = 1833
$ws.Cells.Item(59, 17).Value = 1.5
$ws.Cells.Item(59, 18).Value = "Hortaliza"
